# Insert a new weekly price record as row 276, pushing the existing
# rows 276-307 down to 277-308 (matches the target diff: dimension grows
# from A1:R307 to A1:R308, and every row at/after 276 shifts down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 276; Excel shifts 276:307 -> 277:308 and copies
# the row-above's formatting onto the new row (keeps D's date style s="2").
$ws.Rows("276:276").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(276, 1).Value  = 7
$ws.Cells.Item(276, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(276, 3).Value  = "Ñuble"
$ws.Cells.Item(276, 4).Value  = 44769
$ws.Cells.Item(276, 5).Value  = 16
$ws.Cells.Item(276, 6).Value  = 100114013
$ws.Cells.Item(276, 7).Value  = "Zanahoria"
$ws.Cells.Item(276, 8).Value  = "Sin especificar"
$ws.Cells.Item(276, 9).Value  = "Primera"
$ws.Cells.Item(276, 10).Value = 100
$ws.Cells.Item(276, 11).Value = 9000
$ws.Cells.Item(276, 12).Value = 10000
$ws.Cells.Item(276, 13).Value = 9500
$ws.Cells.Item(276, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(276, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(276, 16).Value = 475
$ws.Cells.Item(276, 17).Value = 20
$ws.Cells.Item(276, 18).Value = "Hortaliza"
